$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Darius Garland, PG, Cleveland Cavaliers -> Austin Reaves, PG,SG, Los Angeles Lakers
$ws.Range("A2").Value = "Austin Reaves"
$ws.Range("B2").Value = "PG,SG"
$ws.Range("C2").Value = "Los Angeles Lakers"

# Row 6: Trey Murphy III, SF,PF, New Orleans Pelicans -> Franz Wagner, SF,PF, Orlando Magic
$ws.Range("A6").Value = "Franz Wagner"
$ws.Range("C6").Value = "Orlando Magic"

# Row 10: Jarrett Allen, C, Cleveland Cavaliers -> Daniel Gafford, PF,C, Dallas Mavericks
$ws.Range("A10").Value = "Daniel Gafford"
$ws.Range("B10").Value = "PF,C"
$ws.Range("C10").Value = "Dallas Mavericks"

# Row 11: Daniel Gafford, PF,C, Dallas Mavericks -> Jalen Duren, C, Detroit Pistons
$ws.Range("A11").Value = "Jalen Duren"
$ws.Range("B11").Value = "C"
$ws.Range("C11").Value = "Detroit Pistons"

# Row 12: Jalen Duren, C, Detroit Pistons -> Tyrese Maxey, PG,SG, Philadelphia 76ers
$ws.Range("A12").Value = "Tyrese Maxey"
$ws.Range("B12").Value = "PG,SG"
$ws.Range("C12").Value = "Philadelphia 76ers"

# Row 13: Franz Wagner, SF,PF, Orlando Magic -> Darius Garland, PG, Cleveland Cavaliers
$ws.Range("A13").Value = "Darius Garland"
$ws.Range("B13").Value = "PG"
$ws.Range("C13").Value = "Cleveland Cavaliers"

# Row 14: Austin Reaves, PG,SG, Los Angeles Lakers -> Ty Jerome, PG,SG, Cleveland Cavaliers
$ws.Range("A14").Value = "Ty Jerome"
$ws.Range("C14").Value = "Cleveland Cavaliers"

# Row 15: Zach Collins, PF,C, San Antonio Spurs -> Jarrett Allen, C, Cleveland Cavaliers
$ws.Range("A15").Value = "Jarrett Allen"
$ws.Range("B15").Value = "C"
$ws.Range("C15").Value = "Cleveland Cavaliers"

# Row 19: Tyrese Maxey, PG,SG, Philadelphia 76ers -> Trey Murphy III, SF,PF, New Orleans Pelicans
$ws.Range("A19").Value = "Trey Murphy III"
$ws.Range("B19").Value = "SF,PF"
$ws.Range("C19").Value = "New Orleans Pelicans"
